$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three rows that are punctuation-variant duplicates of an
# already-present underscore-named row (kept row retained, dotted row
# removed). Delete from the bottom up so row numbers of earlier rows
# stay valid while deleting.
# Original rows (1-indexed):
#  2  patient.id | patient   -> remove (duplicate of row 3 patient_id|patient)
#  5  sample.id  | sample    -> remove (duplicate of row 6 sample_id|sample)
# 10  biopsy.id  | biopsy    -> remove (duplicate of row 9 biopsy_id|biopsy)
$ws.Rows.Item(10).Delete() | Out-Null
$ws.Rows.Item(5).Delete() | Out-Null
$ws.Rows.Item(2).Delete() | Out-Null

# After the three deletions, the table has 11 rows (1 header/general +
# 10 remaining); append a new row with the new pairing.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).Value = "cell_subtype"
$ws.Cells.Item($newRow, 2).Value = "annotation_authors_minor"
